# Auto-generated edit script: updates cached profit calculation values
# across multiple worksheets (ALC, ARM, BSM, CRP, CUL, LTW, WVR) in the
# "Tonberry_Profits" workbook, per scheduled runner refresh.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 806.5574
$ws.Range("I15").Value = 806.5574
$ws.Range("K15").Value = 2419.6722
$ws.Range("M15").Value = -2250.6722
$ws.Range("H32").Value = 1498.0834
$ws.Range("J32").Value = 1523.75
$ws.Range("L32").Value = 1523.75
$ws.Range("N32").Value = -2175.75
$ws.Range("H98").Value = 1032.6471
$ws.Range("I98").Value = 763.625
$ws.Range("J98").Value = 1678.3
$ws.Range("K98").Value = 763.625
$ws.Range("L98").Value = 1678.3
$ws.Range("M98").Value = 734.375
$ws.Range("N98").Value = -4674.3
$ws.Range("H116").Value = 9744.4375
$ws.Range("I116").Value = 35233.332
$ws.Range("J116").Value = 3862.3845
$ws.Range("K116").Value = 35233.332
$ws.Range("L116").Value = 3862.3845
$ws.Range("M116").Value = -31791.332
$ws.Range("N116").Value = -10746.3845
$ws.Range("H122").Value = 1032.6471
$ws.Range("I122").Value = 763.625
$ws.Range("J122").Value = 1678.3
$ws.Range("K122").Value = 2290.875
$ws.Range("L122").Value = 5034.9
$ws.Range("M122").Value = 159.125
$ws.Range("N122").Value = -9934.9
$ws.Range("H137").Value = 1616.5
$ws.Range("I137").Value = 1378.2
$ws.Range("K137").Value = 4134.6
$ws.Range("M137").Value = -1584.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1012.82355
$ws.Range("I74").Value = 509
$ws.Range("K74").Value = 509
$ws.Range("M74").Value = 365
$ws.Range("H77").Value = 1012.82355
$ws.Range("I77").Value = 509
$ws.Range("K77").Value = 2545
$ws.Range("M77").Value = 1823
$ws.Range("H97").Value = 934.6
$ws.Range("I97").Value = 883
$ws.Range("J97").Value = 1399
$ws.Range("K97").Value = 883
$ws.Range("L97").Value = 1399
$ws.Range("M97").Value = -387
$ws.Range("N97").Value = -2391

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1053
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 1053
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 1053
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -4049

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1311.0625
$ws.Range("I22").Value = 466.66666
$ws.Range("J22").Value = 1505.9231
$ws.Range("K22").Value = 466.66666
$ws.Range("L22").Value = 1505.9231
$ws.Range("M22").Value = -116.66666
$ws.Range("N22").Value = -2205.9231
$ws.Range("H31").Value = 3725.2856
$ws.Range("I31").Value = 5460.3335
$ws.Range("J31").Value = 2424
$ws.Range("K31").Value = 5460.3335
$ws.Range("L31").Value = 2424
$ws.Range("M31").Value = -5165.3335
$ws.Range("N31").Value = -3014
$ws.Range("H34").Value = 3725.2856
$ws.Range("I34").Value = 5460.3335
$ws.Range("J34").Value = 2424
$ws.Range("K34").Value = 5460.3335
$ws.Range("L34").Value = 2424
$ws.Range("M34").Value = -5258.3335
$ws.Range("N34").Value = -2828
$ws.Range("H58").Value = 3260.25
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 3260.25
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 3260.25
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -3666.25
$ws.Range("H107").Value = 878.5
$ws.Range("I107").Value = 561.53845
$ws.Range("J107").Value = 4999
$ws.Range("K107").Value = 561.53845
$ws.Range("L107").Value = 4999
$ws.Range("M107").Value = 1358.46155
$ws.Range("N107").Value = -8839
$ws.Range("H134").Value = 1360.68
$ws.Range("I134").Value = 1309.4348
$ws.Range("J134").Value = 1950
$ws.Range("K134").Value = 3928.3044
$ws.Range("L134").Value = 5850
$ws.Range("M134").Value = -1393.3044
$ws.Range("N134").Value = -10920
$ws.Range("H136").Value = 3260.25
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 3260.25
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 9780.75
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -14880.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 14500.5
$ws.Range("J17").Value = 14500.5
$ws.Range("L17").Value = 43501.5
$ws.Range("N17").Value = -43839.5
$ws.Range("H113").Value = 7409.8
$ws.Range("J113").Value = 780.46155
$ws.Range("L113").Value = 2341.38465
$ws.Range("N113").Value = -6681.38465
$ws.Range("H131").Value = 26609.535
$ws.Range("J131").Value = 29720.28
$ws.Range("L131").Value = 89160.84
$ws.Range("N131").Value = -99240.84

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2289.2354
$ws.Range("I61").Value = 1901.4166
$ws.Range("K61").Value = 1901.4166
$ws.Range("M61").Value = -1699.4166
$ws.Range("H68").Value = 1819.125
$ws.Range("I68").Value = 1364.7142
$ws.Range("K68").Value = 1364.7142
$ws.Range("M68").Value = -615.7141999999999
$ws.Range("H71").Value = 1819.125
$ws.Range("I71").Value = 1364.7142
$ws.Range("K71").Value = 6823.571
$ws.Range("M71").Value = -3079.571
$ws.Range("H82").Value = 1969.5834
$ws.Range("I82").Value = 1351.6666
$ws.Range("J82").Value = 3823.3333
$ws.Range("K82").Value = 1351.6666
$ws.Range("L82").Value = 3823.3333
$ws.Range("M82").Value = -990.6666
$ws.Range("N82").Value = -4545.3333
$ws.Range("H85").Value = 1969.5834
$ws.Range("I85").Value = 1351.6666
$ws.Range("J85").Value = 3823.3333
$ws.Range("K85").Value = 1351.6666
$ws.Range("L85").Value = 3823.3333
$ws.Range("M85").Value = -103.6666
$ws.Range("N85").Value = -6319.3333
$ws.Range("H93").Value = 15873960
$ws.Range("I93").Value = 830.8823
$ws.Range("J93").Value = 83334760
$ws.Range("K93").Value = 830.8823
$ws.Range("L93").Value = 83334760
$ws.Range("M93").Value = 417.1177
$ws.Range("N93").Value = -83337256
$ws.Range("H100").Value = 1386.6
$ws.Range("I100").Value = 1363.75
$ws.Range("J100").Value = 1478
$ws.Range("K100").Value = 1363.75
$ws.Range("L100").Value = 1478
$ws.Range("M100").Value = -822.75
$ws.Range("N100").Value = -2560
$ws.Range("H113").Value = 2289.2354
$ws.Range("I113").Value = 1901.4166
$ws.Range("K113").Value = 1901.4166
$ws.Range("M113").Value = 268.5834
$ws.Range("H122").Value = 7421.2856
$ws.Range("I122").Value = 6877.5557
$ws.Range("K122").Value = 20632.6671
$ws.Range("M122").Value = -18182.6671

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1919.4
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 1919.4
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 1919.4
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -4665.4
$ws.Range("H107").Value = 783.82355
$ws.Range("I107").Value = 433.72726
$ws.Range("J107").Value = 1425.6666
$ws.Range("K107").Value = 1301.18178
$ws.Range("L107").Value = 4276.9998
$ws.Range("M107").Value = 618.8182200000001
$ws.Range("N107").Value = -8116.9998
$ws.Range("H136").Value = 2124.4243
$ws.Range("I136").Value = 2234.8333
$ws.Range("J136").Value = 1991.9333
$ws.Range("K136").Value = 6704.499899999999
$ws.Range("L136").Value = 5975.7999
$ws.Range("M136").Value = -4154.499899999999
$ws.Range("N136").Value = -11075.7999
